$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022224115596433
$ws.Range("D2").Value = 1.033587283966038
$ws.Range("E2").Value = 1.0229929622304
$ws.Range("F2").Value = 1.020668514257566
$ws.Range("I2").Value = 1.032842740199678
$ws.Range("J2").Value = 1.027411156560433
$ws.Range("K2").Value = 1.036389204040683
$ws.Range("L2").Value = 1.025825709666696
$ws.Range("M2").Value = 1.023508118530196
$ws.Range("N2").Value = 1.028870198248457

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023208443316025
$ws.Range("D3").Value = 1.034111848828131
$ws.Range("E3").Value = 1.023829291529974
$ws.Range("F3").Value = 1.022291480738175
$ws.Range("I3").Value = 1.033014980912342
$ws.Range("J3").Value = 1.028033260391826
$ws.Range("K3").Value = 1.03672353958252
$ws.Range("L3").Value = 1.026468742571933
$ws.Range("M3").Value = 1.024935135452937
$ws.Range("N3").Value = 1.029493185538646

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023845235939103
$ws.Range("D4").Value = 1.034451234842616
$ws.Range("E4").Value = 1.024370713133348
$ws.Range("F4").Value = 1.02334140932584
$ws.Range("I4").Value = 1.033125176898019
$ws.Range("J4").Value = 1.028435100753629
$ws.Range("K4").Value = 1.036939122528783
$ws.Range("L4").Value = 1.026884447277201
$ws.Range("M4").Value = 1.025857813523343
$ws.Range("N4").Value = 1.029895596559854

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.024112912389566
$ws.Range("D5").Value = 1.034593901238995
$ws.Range("E5").Value = 1.024598389197918
$ws.Range("F5").Value = 1.023782748463186
$ws.Range("I5").Value = 1.033171202562052
$ws.Range("J5").Value = 1.02860386673298
$ws.Range("K5").Value = 1.037029572221397
$ws.Range("L5").Value = 1.02705911866784
$ws.Range("M5").Value = 1.02624554599183
$ws.Range("N5").Value = 1.030064602206254

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.024157854636789
$ws.Range("D6").Value = 1.034617854843222
$ws.Range("E6").Value = 1.024636620672185
$ws.Range("F6").Value = 1.02385684845175
$ws.Range("I6").Value = 1.033178912832466
$ws.Range("J6").Value = 1.02863219346201
$ws.Range("K6").Value = 1.037044748464843
$ws.Range("L6").Value = 1.027088441441421
$ws.Range("M6").Value = 1.026310638676097
$ws.Range("N6").Value = 1.030092969162487

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023848812766505
$ws.Range("D7").Value = 1.034453141204769
$ws.Range("E7").Value = 1.02437375510677
$ws.Range("F7").Value = 1.023347306712772
$ws.Range("I7").Value = 1.033125793077124
$ws.Range("J7").Value = 1.028437356471338
$ws.Range("K7").Value = 1.036940331835646
$ws.Range("L7").Value = 1.026886781601628
$ws.Range("M7").Value = 1.02586299505085
$ws.Range("N7").Value = 1.029897855480941

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022556801783996
$ws.Range("D8").Value = 1.033764571142258
$ws.Range("E8").Value = 1.023275549689091
$ws.Range("F8").Value = 1.021217056794785
$ws.Range("I8").Value = 1.03290120953988
$ws.Range("J8").Value = 1.027621545176213
$ws.Range("K8").Value = 1.036502350079921
$ws.Range("L8").Value = 1.026043104599276
$ws.Range("M8").Value = 1.023990533301224
$ws.Range("N8").Value = 1.029080885640202

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020279059721823
$ws.Range("D9").Value = 1.032550954106086
$ws.Range("E9").Value = 1.021342365328781
$ws.Range("F9").Value = 1.017461163589779
$ws.Range("I9").Value = 1.03249585991647
$ws.Range("J9").Value = 1.026178580956326
$ws.Range("K9").Value = 1.035724824318101
$ws.Range("L9").Value = 1.024553511478048
$ws.Range("M9").Value = 1.020685427601356
$ws.Range("N9").Value = 1.027635872245612

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018759807181394
$ws.Range("D10").Value = 1.031741771602007
$ws.Range("E10").Value = 1.020054910167736
$ws.Range("F10").Value = 1.014955404961074
$ws.Range("I10").Value = 1.032219180474266
$ws.Range("J10").Value = 1.025212942370349
$ws.Range("K10").Value = 1.03520265363465
$ws.Range("L10").Value = 1.02355846304764
$ws.Range("M10").Value = 1.018477917812245
$ws.Range("N10").Value = 1.026668862342083

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018101762388722
$ws.Range("D11").Value = 1.031391377649892
$ws.Range("E11").Value = 1.019497742024176
$ws.Range("F11").Value = 1.013869859628069
$ws.Range("I11").Value = 1.032097848274942
$ws.Range("J11").Value = 1.024793933781361
$ws.Range("K11").Value = 1.034975650642298
$ws.Range("L11").Value = 1.023127120000327
$ws.Range("M11").Value = 1.017520989961643
$ws.Range("N11").Value = 1.026249258712838

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.0178573043098
$ws.Range("D12").Value = 1.031261225090957
$ws.Range("E12").Value = 1.019290831003438
$ws.Range("F12").Value = 1.013466551511058
$ws.Range("I12").Value = 1.032052550615211
$ws.Range("J12").Value = 1.024638162263423
$ws.Range("K12").Value = 1.034891197129424
$ws.Range("L12").Value = 1.022966827280753
$ws.Range("M12").Value = 1.01716537809556
$ws.Range("N12").Value = 1.026093265981477

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017909742816114
$ws.Range("D13").Value = 1.031289143265882
$ws.Range("E13").Value = 1.019335212041534
$ws.Range("F13").Value = 1.01355306659055
$ws.Range("I13").Value = 1.032062277499606
$ws.Range("J13").Value = 1.024671581816788
$ws.Range("K13").Value = 1.034909318770136
$ws.Range("L13").Value = 1.0230012139117
$ws.Range("M13").Value = 1.017241665680432
$ws.Range("N13").Value = 1.026126732994441

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018081556045484
$ws.Range("D14").Value = 1.031380619204729
$ws.Range("E14").Value = 1.019480637757335
$ws.Range("F14").Value = 1.013836523898886
$ws.Range("I14").Value = 1.032094108636912
$ws.Range("J14").Value = 1.024781060378115
$ws.Range("K14").Value = 1.03496867242839
$ws.Range("L14").Value = 1.023113871640319
$ws.Range("M14").Value = 1.017491598404691
$ws.Range("N14").Value = 1.026236367027883

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018187411712285
$ws.Range("D15").Value = 1.03143698050155
$ws.Range("E15").Value = 1.019570245431414
$ws.Range("F15").Value = 1.014011159373672
$ws.Range("I15").Value = 1.032113690446339
$ws.Range("J15").Value = 1.024848496076979
$ws.Range("K15").Value = 1.03500522437419
$ws.Range("L15").Value = 1.023183274135257
$ws.Range("M15").Value = 1.01764556797142
$ws.Range("N15").Value = 1.026303898493173

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018803475170087
$ws.Range("D16").Value = 1.031765025936122
$ws.Range("E16").Value = 1.020091894055182
$ws.Range("F16").Value = 1.015027437084492
$ws.Range("I16").Value = 1.032227200705511
$ws.Range("J16").Value = 1.025240731970653
$ws.Range("K16").Value = 1.035217700178304
$ws.Range("L16").Value = 1.023587079697535
$ws.Range("M16").Value = 1.018541403019481
$ws.Range("N16").Value = 1.026696691406807

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.01918986158204
$ws.Range("D17").Value = 1.031970797684243
$ws.Range("E17").Value = 1.020419192821765
$ws.Range("F17").Value = 1.015664772416655
$ws.Range("I17").Value = 1.032297993636738
$ws.Range("J17").Value = 1.02548653486328
$ws.Range("K17").Value = 1.03535074025021
$ws.Range("L17").Value = 1.023840247251533
$ws.Range("M17").Value = 1.019103047704455
$ws.Range("N17").Value = 1.026942843367736

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019415215171143
$ws.Range("D18").Value = 1.032090819565912
$ws.Range("E18").Value = 1.0206101304355
$ws.Range("F18").Value = 1.01603646820991
$ws.Range("I18").Value = 1.032339138501692
$ws.Range("J18").Value = 1.025629822607367
$ws.Range("K18").Value = 1.035428253426692
$ws.Range("L18").Value = 1.023987869275241
$ws.Range("M18").Value = 1.019430543579059
$ws.Range("N18").Value = 1.027086334596855

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.019492051721475
$ws.Range("D19").Value = 1.032131743667707
$ws.Range("E19").Value = 1.020675240287899
$ws.Range("F19").Value = 1.01616319843595
$ws.Range("I19").Value = 1.032353142823122
$ws.Range("J19").Value = 1.025678665643379
$ws.Range("K19").Value = 1.035454668665678
$ws.Range("L19").Value = 1.024038196728051
$ws.Range("M19").Value = 1.019542194119004
$ws.Range("N19").Value = 1.027135246995579

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019148407956999
$ws.Range("D20").Value = 1.031948720448233
$ws.Range("E20").Value = 1.02038407368331
$ws.Range("F20").Value = 1.015596397710755
$ws.Range("I20").Value = 1.032290413476202
$ws.Range("J20").Value = 1.025460171328577
$ws.Range("K20").Value = 1.035336475287657
$ws.Range("L20").Value = 1.023813089564352
$ws.Range("M20").Value = 1.019042799140277
$ws.Range("N20").Value = 1.026916442393789

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.018030962195104
$ws.Range("D21").Value = 1.031353681828524
$ws.Range("E21").Value = 1.019437812240349
$ws.Range("F21").Value = 1.013753055285608
$ws.Range("I21").Value = 1.03208474149419
$ws.Range("J21").Value = 1.024748825350021
$ws.Range("K21").Value = 1.034951197962764
$ws.Range("L21").Value = 1.023080698758387
$ws.Range("M21").Value = 1.017418004090795
$ws.Range("N21").Value = 1.026204086222351

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017328200096339
$ws.Range("D22").Value = 1.030979554315895
$ws.Range("E22").Value = 1.01884312630464
$ws.Range("F22").Value = 1.012593558273365
$ws.Range("I22").Value = 1.031954099717142
$ws.Range("J22").Value = 1.024300803608395
$ws.Range("K22").Value = 1.034708180999909
$ws.Range("L22").Value = 1.022619794672838
$ws.Range("M22").Value = 1.016395465678021
$ws.Range("N22").Value = 1.025755428238464

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017700764986282
$ws.Range("D23").Value = 1.031177886270729
$ws.Range("E23").Value = 1.019158355443394
$ws.Range("F23").Value = 1.013208280966383
$ws.Range("I23").Value = 1.032023481203083
$ws.Range("J23").Value = 1.024538381679816
$ws.Range("K23").Value = 1.034837082397292
$ws.Range("L23").Value = 1.022864168721672
$ws.Range("M23").Value = 1.016937626404931
$ws.Range("N23").Value = 1.025993343697997

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019167139124349
$ws.Range("D24").Value = 1.031958696205007
$ws.Range("E24").Value = 1.020399942418076
$ws.Range("F24").Value = 1.01562729345335
$ws.Range("I24").Value = 1.032293839080031
$ws.Range("J24").Value = 1.025472084137233
$ws.Range("K24").Value = 1.035342921278581
$ws.Range("L24").Value = 1.02382536109712
$ws.Range("M24").Value = 1.019070023186533
$ws.Range("N24").Value = 1.02692837212

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020868040564191
$ws.Range("D25").Value = 1.032864727150826
$ws.Range("E25").Value = 1.021841904135908
$ws.Range("F25").Value = 1.018432444818276
$ws.Range("I25").Value = 1.032601789543612
$ws.Range("J25").Value = 1.026552264989535
$ws.Range("K25").Value = 1.035926509826178
$ws.Range("L25").Value = 1.024938955706894
$ws.Range("M25").Value = 1.021540574735618
$ws.Range("N25").Value = 1.02801008695301
